# "Added filter button, working on adding integrated db"
#
# - Sheet1: selection nudges from F9 to G9 (room made on the sheet for the
#   new filter button).
# - New "Sheet2" is added after "Sheet1" to hold the column headers for the
#   integrated card database (uuid / name / variation json / detail json /
#   image blob) and becomes the active tab.

$wb = $excel.ActiveWorkbook

# --- Sheet1: move the selection one column over ----------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
[void]$ws1.Range("G9").Select()

# --- Sheet2: new sheet for the integrated DB table headers -----------------
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "uuid"
$ws2.Range("B1").Value = "name"
$ws2.Range("C1").Value = "variation json"
$ws2.Range("D1").Value = "detail json"
$ws2.Range("E1").Value = "image blob"

$ws2.Columns.Item(3).ColumnWidth = 14.109375
$ws2.Columns.Item(4).ColumnWidth = 11.5546875
$ws2.Columns.Item(5).ColumnWidth = 11.6640625

[void]$ws2.Range("A2").Select()
